# Add new corona-cases rows (171-179) for date 2020-03-27 (serial 43917)
# to the KSA corona cases worksheet, extending the running-total formula
# in column D and updating the view selection, per the commit:
# "added visualizations, updated ksa dataset"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: dates (same date format as the rest of the column) -------
$ws.Range("A171:A179").Value = 43917
$ws.Range("A171:A179").NumberFormat = "yyyy\-mm\-dd"

# --- Column B: city ------------------------------------------------------
$ws.Range("B171").Value = "Riyadh"
$ws.Range("B172").Value = "Medinah"
$ws.Range("B173").Value = "Qatif"
$ws.Range("B174").Value = "Jeddah"
$ws.Range("B175").Value = "Dammam"
$ws.Range("B176").Value = "Dhahran"
$ws.Range("B177").Value = "Buraidah"
$ws.Range("B178").Value = "Ahsaa"
$ws.Range("B179").Value = "Kobar"

# --- Column C: province ---------------------------------------------------
$ws.Range("C171").Value = "Riyadh"
$ws.Range("C172").Value = "Medinah"
$ws.Range("C173").Value = "Eastern province"
$ws.Range("C174").Value = "Makkah"
$ws.Range("C175").Value = "Eastern province"
$ws.Range("C176").Value = "Eastern province"
$ws.Range("C177").Value = "Al Quassim"
$ws.Range("C178").Value = "Eastern province"
$ws.Range("C179").Value = "Eastern province"

# --- Column D: running total (same relative formula as the rest of the
#     column: each row = previous row's total + this row's new_cases) -----
$ws.Range("D171:D179").Formula = "=D170+E171"

# --- Column E: new cases ---------------------------------------------------
$ws.Range("E171").Value = 46
$ws.Range("E172").Value = 19
$ws.Range("E173").Value = 10
$ws.Range("E174").Value = 7
$ws.Range("E175").Value = 4
$ws.Range("E176").Value = 2
$ws.Range("E177").Value = 2
$ws.Range("E178").Value = 1
$ws.Range("E179").Value = 1

# --- Column H: map_name -----------------------------------------------------
$ws.Range("H171").Value = "Ar Riyad"
$ws.Range("H172").Value = "Al Madinah"
$ws.Range("H173").Value = "Ash Sharqiyah"
$ws.Range("H174").Value = "Makkah"
$ws.Range("H175").Value = "Ash Sharqiyah"
$ws.Range("H176").Value = "Ash Sharqiyah"
$ws.Range("H177").Value = "Al Quassim"
$ws.Range("H178").Value = "Ash Sharqiyah"
$ws.Range("H179").Value = "Ash Sharqiyah"

# --- Update the view: scroll down and select the next empty cell ---------
$ws.Range("A151").Select()
$ws.Range("H180").Select()
